$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 4.2
$ws.Range("G2").Value = 6
$ws.Range("H2").Value = 1.74
$ws.Range("I2").Value = 1.93
$ws.Range("K2").Value = 4.8
$ws.Range("N2").Value = 4.1
$ws.Range("O2").Value = 1.27
$ws.Range("P2").Value = 1.94
$ws.Range("Q2").Value = 1.71
$ws.Range("R2").Value = 1.42
$ws.Range("S2").Value = 2.8
$ws.Range("U2").Value = 2.1
$ws.Range("V2").Value = 2.06
$ws.Range("W2").Value = 1.2
$ws.Range("Y2").Value = 11
$ws.Range("AO2").Value = 12

# Row 3
$ws.Range("F3").Value = 2.68
$ws.Range("G3").Value = 2.9
$ws.Range("H3").Value = 3.25
$ws.Range("I3").Value = 3.55
$ws.Range("J3").Value = 2.76
$ws.Range("K3").Value = 3
$ws.Range("N3").Value = 2.26
$ws.Range("O3").Value = 1.7
$ws.Range("Q3").Value = 3.15
$ws.Range("T3").Value = 2.34
$ws.Range("V3").Value = 1.39
$ws.Range("AK3").Value = 1000

# Row 4
$ws.Range("I4").Value = 6.4
$ws.Range("L4").Value = 1.33
$ws.Range("N4").Value = 3.55
$ws.Range("O4").Value = 1.27
$ws.Range("Q4").Value = 1.84
$ws.Range("AF4").Value = 12.5

# Row 5
$ws.Range("F5").Value = 1.67
$ws.Range("J5").Value = 3.7
$ws.Range("Q5").Value = 2.06
$ws.Range("T5").Value = 1.98

# Row 6
$ws.Range("H6").Value = 15
$ws.Range("I6").Value = 18
$ws.Range("K6").Value = 8.6
$ws.Range("L6").Value = 1.2
$ws.Range("N6").Value = 6.6
$ws.Range("O6").Value = 1.14
$ws.Range("P6").Value = 2.9
$ws.Range("Q6").Value = 1.41
$ws.Range("R6").Value = 1.76
$ws.Range("S6").Value = 2.06
$ws.Range("T6").Value = 2.02
$ws.Range("U6").Value = 1.79
$ws.Range("X6").Value = 42
$ws.Range("Y6").Value = 70
$ws.Range("Z6").Value = 170
$ws.Range("AB6").Value = 12.5
$ws.Range("AD6").Value = 1000
$ws.Range("AE6").Value = 290
$ws.Range("AH6").Value = 38
$ws.Range("AI6").Value = 170
$ws.Range("AJ6").Value = 11
$ws.Range("AL6").Value = 40
$ws.Range("AN6").Value = 3.45

# Row 7
$ws.Range("J7").Value = 7.8
$ws.Range("K7").Value = 11
$ws.Range("N7").Value = 6.4
$ws.Range("O7").Value = 1.15
$ws.Range("P7").Value = 2.92
$ws.Range("Q7").Value = 1.43
$ws.Range("R7").Value = 1.75
$ws.Range("T7").Value = 2.22
$ws.Range("V7").Value = 5.7
$ws.Range("Y7").Value = 13
$ws.Range("Z7").Value = 9.4
$ws.Range("AA7").Value = 9.800000000000001

# Row 8
$ws.Range("F8").Value = 2.84
$ws.Range("G8").Value = 3
$ws.Range("H8").Value = 3.2
$ws.Range("I8").Value = 3.45
$ws.Range("J8").Value = 2.72
$ws.Range("K8").Value = 2.94
$ws.Range("M8").Value = 1.18
$ws.Range("N8").Value = 2.18
$ws.Range("O8").Value = 1.75
$ws.Range("U8").Value = 1.62
$ws.Range("V8").Value = 1.4
$ws.Range("W8").Value = 1.5
$ws.Range("Y8").Value = 7.8
$ws.Range("AA8").Value = 90
$ws.Range("AB8").Value = 7.4
$ws.Range("AI8").Value = 130
$ws.Range("AO8").Value = 120

# Row 9
$ws.Range("N9").Value = 2.76
$ws.Range("O9").Value = 1.5
$ws.Range("Q9").Value = 2.46
$ws.Range("S9").Value = 4.3
$ws.Range("T9").Value = 2.22
$ws.Range("U9").Value = 1.7
$ws.Range("AB9").Value = 7.4
$ws.Range("AC9").Value = 9.800000000000001

# Row 10
$ws.Range("F10").Value = 1.73
$ws.Range("O10").Value = 1.4
$ws.Range("P10").Value = 1.7
$ws.Range("Q10").Value = 2.2
$ws.Range("T10").Value = 2.02
$ws.Range("AL10").Value = 50

# Row 11
$ws.Range("F11").Value = 2.52
$ws.Range("G11").Value = 2.86
$ws.Range("H11").Value = 3.15
$ws.Range("J11").Value = 2.82
$ws.Range("K11").Value = 3.25
$ws.Range("N11").Value = 2.44
$ws.Range("U11").Value = 1.73
$ws.Range("Y11").Value = 11.5
$ws.Range("AB11").Value = 9.199999999999999
$ws.Range("AE11").Value = 70

# Row 12
$ws.Range("G12").Value = 2.34
$ws.Range("H12").Value = 3.6
$ws.Range("I12").Value = 3.8
$ws.Range("S12").Value = 3.8
$ws.Range("W12").Value = 1.74

# Row 13
$ws.Range("G13").Value = 2.64
$ws.Range("I13").Value = 4.5
$ws.Range("J13").Value = 2.74
$ws.Range("N13").Value = 2.8
$ws.Range("O13").Value = 1.45
$ws.Range("P13").Value = 1.61
$ws.Range("Q13").Value = 2.32
$ws.Range("R13").Value = 1.22
$ws.Range("S13").Value = 4.1
$ws.Range("W13").Value = 1.61
$ws.Range("X13").Value = 12
$ws.Range("Y13").Value = 14
$ws.Range("AC13").Value = 8.800000000000001
$ws.Range("AE13").Value = 75
$ws.Range("AI13").Value = 95

# Row 14
$ws.Range("G14").Value = 1.7
$ws.Range("O14").Value = 1.35
$ws.Range("T14").Value = 2
$ws.Range("W14").Value = 2.42
$ws.Range("AB14").Value = 7.6
$ws.Range("AD14").Value = 25
$ws.Range("AO14").Value = 140

# Row 15
$ws.Range("G15").Value = 1.51
$ws.Range("I15").Value = 10
$ws.Range("Z15").Value = 100

# Row 16
$ws.Range("Q16").Value = 2.42
$ws.Range("R16").Value = 1.22
$ws.Range("AE16").Value = 310
$ws.Range("AI16").Value = 280
$ws.Range("AM16").Value = 410

Write-Host "Applied 149 cell updates"
